$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New header cells for columns I and J, matching the header style used by H1
$ws.Range("I1").Value = "I0"
$ws.Range("J1").Value = "IF"

$ws.Range("H1").Copy() | Out-Null
$ws.Range("I1:J1").PasteSpecial(-4122) | Out-Null

# New data values for columns I and J
$ws.Range("I2").Value = 7
$ws.Range("J2").Value = 7

$ws.Range("I3").Value = 9
$ws.Range("J3").Value = 9

$ws.Range("I4").Value = 8
$ws.Range("J4").Value = 8
